$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collection Date (text, not a date value)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "24-04-2024"

# Latitude / Longitude
$ws.Range("C2").Value = 654
$ws.Range("D2").Value = 564

# Name
$ws.Range("E2").Value = "asdjkhajksdh"

# Area (ha)
$ws.Range("F2").Value = 2.5

# Gender
$ws.Range("G2").Value = "Female"

# Age
$ws.Range("H2").Value = 54

# Address
$ws.Range("I2").Value = "654asdasd"

# Mobile No.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "5646545645"

# Soil pH
$ws.Range("K2").Value = 2

# Nitrogen
$ws.Range("L2").Value = 56

# Phosphorus
$ws.Range("M2").Value = 48

# Potassium
$ws.Range("N2").Value = 45

# Electrical Conductivity
$ws.Range("O2").Value = 1

# Temperature
$ws.Range("P2").Value = 5

# Moisture
$ws.Range("Q2").Value = 23

# Humidity
$ws.Range("R2").Value = 21

# Soil Health Score
$ws.Range("S2").Value = 0.1714898508159927

# Recommendations
$ws.Range("T2").Value = "Focus on soil remediation. Grow leguminous cover crops like cowpea, horse gram, or sunn hemp."

# Fertilizer Recommendation
$ws.Range("U2").Value = "Apply organic amendments like compost (5-7.5 tonnes/ha), vermicompost (2.5-3.5 tonnes/ha), or well-decomposed farmyard manure (10-12.5 tonnes/ha). Incorporate green manure crops like dhaincha (Sesbania aculeata) (5-6 tonnes/ha), sunhemp (Crotalaria juncea) (4-5 tonnes/ha), or cowpea (Vigna unguiculata) (3-4 tonnes/ha). Avoid applying chemical fertilizers until soil health improves."
